$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1861.4166
$ws.Range("I28").Value = 158
$ws.Range("K28").Value = 158
$ws.Range("M28").Value = 327
$ws.Range("H62").Value = 6180.1763
$ws.Range("I62").Value = 4024.75
$ws.Range("K62").Value = 4024.75
$ws.Range("M62").Value = -3400.75
$ws.Range("H65").Value = 6180.1763
$ws.Range("I65").Value = 4024.75
$ws.Range("K65").Value = 20123.75
$ws.Range("M65").Value = -17003.75
$ws.Range("H74").Value = 7248.8213
$ws.Range("J74").Value = 7682.72
$ws.Range("L74").Value = 7682.72
$ws.Range("N74").Value = -9554.720000000001
$ws.Range("H77").Value = 7248.8213
$ws.Range("J77").Value = 7682.72
$ws.Range("L77").Value = 38413.6
$ws.Range("N77").Value = -47773.6
$ws.Range("H80").Value = 5546.7144
$ws.Range("J80").Value = 6445.9414
$ws.Range("L80").Value = 19337.8242
$ws.Range("N80").Value = -21333.8242
$ws.Range("H83").Value = 5546.7144
$ws.Range("J83").Value = 6445.9414
$ws.Range("L83").Value = 58013.47259999999
$ws.Range("N83").Value = -67997.47259999999
$ws.Range("H100").Value = 2421.889
$ws.Range("I100").Value = 2207.4666
$ws.Range("J100").Value = 3494
$ws.Range("K100").Value = 2207.4666
$ws.Range("L100").Value = 3494
$ws.Range("M100").Value = -1666.4666
$ws.Range("N100").Value = -4576
$ws.Range("H127").Value = 7000
$ws.Range("J127").Value = 7000
$ws.Range("L127").Value = 21000
$ws.Range("N127").Value = -30920
$ws.Range("H141").Value = 4973.7144
$ws.Range("I141").Value = 5760.8096
$ws.Range("J141").Value = 2612.4285
$ws.Range("K141").Value = 17282.4288
$ws.Range("L141").Value = 7837.2855
$ws.Range("M141").Value = -12102.4288
$ws.Range("N141").Value = -18197.2855
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 6026
$ws.Range("I3").Value = 3368
$ws.Range("J3").Value = 14000
$ws.Range("K3").Value = 3368
$ws.Range("L3").Value = 14000
$ws.Range("M3").Value = -3253
$ws.Range("N3").Value = -14230
$ws.Range("H44").Value = 6944
$ws.Range("H45").Value = 9596071
$ws.Range("H61").Value = 5422.3213
$ws.Range("I61").Value = 5721.2
$ws.Range("K61").Value = 5721.2
$ws.Range("M61").Value = -5509.2
$ws.Range("H97").Value = 1550304
$ws.Range("I97").Value = 2317478.2
$ws.Range("J97").Value = 15955.714
$ws.Range("K97").Value = 2317478.2
$ws.Range("L97").Value = 15955.714
$ws.Range("M97").Value = -2316982.2
$ws.Range("N97").Value = -16947.714
$ws.Range("H102").Value = 4392614
$ws.Range("I102").Value = 5561911
$ws.Range("J102").Value = 7748.75
$ws.Range("K102").Value = 5561911
$ws.Range("L102").Value = 7748.75
$ws.Range("M102").Value = -5560289
$ws.Range("N102").Value = -10992.75
$ws.Range("H122").Value = 909190.5
$ws.Range("I122").Value = 3511.0527
$ws.Range("K122").Value = 10533.1581
$ws.Range("M122").Value = -8083.158100000001
$ws.Range("H136").Value = 5422.3213
$ws.Range("I136").Value = 5721.2
$ws.Range("K136").Value = 17163.6
$ws.Range("M136").Value = -14613.6
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 867.6667
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H99").Value = 8994043
$ws.Range("I99").Value = 11990408
$ws.Range("J99").Value = 4948.75
$ws.Range("K99").Value = 11990408
$ws.Range("L99").Value = 4948.75
$ws.Range("M99").Value = -11988910
$ws.Range("N99").Value = -7944.75
$ws.Range("H105").Value = 3474560.2
$ws.Range("I105").Value = 3474560.2
$ws.Range("K105").Value = 3474560.2
$ws.Range("M105").Value = -3472813.2
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7881.8076
$ws.Range("I31").Value = 8612.842000000001
$ws.Range("K31").Value = 8612.842000000001
$ws.Range("M31").Value = -8317.842000000001
$ws.Range("H34").Value = 7881.8076
$ws.Range("I34").Value = 8612.842000000001
$ws.Range("K34").Value = 8612.842000000001
$ws.Range("M34").Value = -8410.842000000001
$ws.Range("H99").Value = 4447
$ws.Range("I99").Value = 3243.75
$ws.Range("K99").Value = 3243.75
$ws.Range("M99").Value = -1745.75
$ws.Range("H105").Value = 1237.6364
$ws.Range("I105").Value = 1237.6364
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1237.6364
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 509.3635999999999
$ws.Range("N105").ClearContents()
$ws.Range("H106").Value = 29994.5
$ws.Range("J106").Value = 29994.5
$ws.Range("L106").Value = 29994.5
$ws.Range("N106").Value = -32518.5
$ws.Range("H126").Value = 4447
$ws.Range("I126").Value = 3243.75
$ws.Range("K126").Value = 9731.25
$ws.Range("M126").Value = -7261.25
$ws.Range("H132").Value = 40524.96
$ws.Range("I132").Value = 47498.316
$ws.Range("K132").Value = 142494.948
$ws.Range("M132").Value = -139964.948
$ws.Range("H134").Value = 22795.38
$ws.Range("I134").Value = 27880.975
$ws.Range("K134").Value = 83642.92499999999
$ws.Range("M134").Value = -81107.92499999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 289.6129
$ws.Range("J2").Value = 376.94736
$ws.Range("L2").Value = 2261.68416
$ws.Range("N2").Value = -2487.68416
$ws.Range("H7").Value = 3254.2222
$ws.Range("I7").Value = 5334
$ws.Range("K7").Value = 16002
$ws.Range("M7").Value = -15890
$ws.Range("H23").Value = 392.5
$ws.Range("J23").Value = 392.5
$ws.Range("L23").Value = 1177.5
$ws.Range("N23").Value = -1647.5
$ws.Range("H38").Value = 59.9375
$ws.Range("I38").Value = 26.571428
$ws.Range("J38").Value = 293.5
$ws.Range("K38").Value = 79.71428400000001
$ws.Range("L38").Value = 880.5
$ws.Range("M38").Value = 267.285716
$ws.Range("N38").Value = -1574.5
$ws.Range("H56").Value = 15631307
$ws.Range("I56").Value = 15631307
$ws.Range("K56").Value = 15631307
$ws.Range("M56").Value = -15630777
$ws.Range("H136").Value = 2198.1667
$ws.Range("I136").Value = 2198.1667
$ws.Range("K136").Value = 6594.500100000001
$ws.Range("M136").Value = -1494.500100000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H99").Value = 3498.8
$ws.Range("I99").Value = 3498.8
$ws.Range("K99").Value = 3498.8
$ws.Range("M99").Value = -1252.8
$ws.Range("H113").Value = 8773725
$ws.Range("I113").Value = 15153006
$ws.Range("J113").Value = 2212.25
$ws.Range("K113").Value = 15153006
$ws.Range("L113").Value = 2212.25
$ws.Range("M113").Value = -15150836
$ws.Range("N113").Value = -6552.25
$ws.Range("H126").Value = 4295890.5
$ws.Range("I126").Value = 1686422.1
$ws.Range("J126").Value = 27781106
$ws.Range("K126").Value = 5059266.300000001
$ws.Range("L126").Value = 83343318
$ws.Range("M126").Value = -5056796.300000001
$ws.Range("N126").Value = -83348258
$ws.Range("H132").Value = 9837
$ws.Range("I132").Value = 7134.615
$ws.Range("J132").Value = 27402.5
$ws.Range("K132").Value = 21403.845
$ws.Range("L132").Value = 82207.5
$ws.Range("M132").Value = -18873.845
$ws.Range("N132").Value = -87267.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2069
$ws.Range("I22").Value = 1700
$ws.Range("K22").Value = 1700
$ws.Range("M22").Value = -1405
$ws.Range("H27").Value = 2069
$ws.Range("I27").Value = 1700
$ws.Range("K27").Value = 1700
$ws.Range("M27").Value = -1593
$ws.Range("H46").Value = 1896963.6
$ws.Range("I46").Value = 2722847.5
$ws.Range("J46").Value = 9229
$ws.Range("K46").Value = 2722847.5
$ws.Range("L46").Value = 9229
$ws.Range("M46").Value = -2722659.5
$ws.Range("N46").Value = -9605
$ws.Range("H100").Value = 32015.314
$ws.Range("I100").Value = 3443.7693
$ws.Range("K100").Value = 3443.7693
$ws.Range("M100").Value = -2902.7693
$ws.Range("H101").Value = 14086.4
$ws.Range("J101").Value = 14086.4
$ws.Range("L101").Value = 14086.4
$ws.Range("N101").Value = -20576.4
$ws.Range("H132").Value = 11795.692
$ws.Range("I132").Value = 13117.594
$ws.Range("K132").Value = 39352.782
$ws.Range("M132").Value = -36822.782
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 9999
$ws.Range("J15").Value = 9999
$ws.Range("L15").Value = 9999
$ws.Range("N15").Value = -10575
$ws.Range("H107").Value = 38467624
$ws.Range("I107").Value = 55559064
$ws.Range("K107").Value = 166677192
$ws.Range("M107").Value = -166675272
$ws.Range("H126").Value = 2622.7144
$ws.Range("I126").Value = 2661.5
$ws.Range("J126").Value = 2498.6
$ws.Range("K126").Value = 7984.5
$ws.Range("L126").Value = 7495.799999999999
$ws.Range("M126").Value = -5514.5
$ws.Range("N126").Value = -12435.8
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 40000
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 40000
$ws.Range("M130").ClearContents()
$ws.Range("N130").Value = -50040
$ws.Range("H132").Value = 20209070
$ws.Range("I132").Value = 21746522
$ws.Range("K132").Value = 65239566
$ws.Range("M132").Value = -65237036
